$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "theta_threshold_range" parameter row (row 5) was removed entirely.
# Deleting the row shifts the "pie_threshold_range" row up from row 6 to
# row 5, and drops the now-unused shared string automatically on save.
$ws.Rows(5).Delete()

# Update the remaining threshold values.
$ws.Range("B2").Value = 3.8
$ws.Range("C2").Value = 12.7

$ws.Range("B3").Value = 5.3
$ws.Range("C3").Value = 11

$ws.Range("B4").Value = 0.7
$ws.Range("C4").Value = 1.3

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 20

# Restore the single-cell selection left by the author.
$null = $ws.Range("C3").Select()

# Page setup touched during the edit (paper size / orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
